# "Compil - modif staff" - update staff.xlsx with latest roster info

$wb = $excel.ActiveWorkbook

$wsCO          = $wb.Worksheets.Item("CO")
$wsSOUTIEN     = $wb.Worksheets.Item("SOUTIEN")
$wsCOMM        = $wb.Worksheets.Item("COMM")
$wsCOMM_SPRINT = $wb.Worksheets.Item("COMM_SPRINT")

# --- SOUTIEN: Radio-Tour contact person updated ---
$wsSOUTIEN.Range("C2").Value = "Bruno Gauthier"

# --- COMM: a new row for "Commissaires in training / stagiaires" is
#     inserted right after the "Commissaires" row, and the Commissaires
#     roster picks up a new name instead of a TBD placeholder.
$wsCOMM.Range("C3").Value = "Steve Head<br/>`nAndrew Paradowski<br/>`nHélène Soulard<br/>`nNancy Daigle<br/>`nFélix-Antoine Malo<br/>`nGeneviève Marcotte"

$wsCOMM.Rows("4").Insert()
$wsCOMM.Range("A4").Value = "Commissaires in training"
$wsCOMM.Range("B4").Value = "Commissaires stagiaires"
$wsCOMM.Range("C4").Value = "Julie Barbeau<br/>`nDaniel Caron"
$wsCOMM.Rows("4").RowHeight = 34
$wsCOMM.PageSetup.Orientation = 1

# --- COMM_SPRINT: same roster updates as COMM, plus the "Results and
#     photofinish" formula now points one row further down because of
#     the row inserted above in COMM.
$wsCOMM_SPRINT.Range("C2").Value = "Andrew Paradowski"
$wsCOMM_SPRINT.Range("C3").Value = "Nancy Daigle<br/>`nGeneviève Marcotte"
$wsCOMM_SPRINT.Range("C4").Formula = "=COMM!C5"

# --- Restore each sheet's selection / active cell, then make COMM the
#     active tab, matching where the author ended up after editing.
$wsCO.Range("A27").Select() | Out-Null
$wsSOUTIEN.Range("C21").Select() | Out-Null
$wsCOMM_SPRINT.Range("C10").Select() | Out-Null
$wsCOMM.Range("A5").Select() | Out-Null
$wsCOMM.Activate() | Out-Null
